$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted at row 116 (between the
# existing "La Araucanía" record that was at row 116 and the rest of the
# data). Insert a blank row there; this pushes the old rows 116-215 down to
# 117-216 (so the prior last row, 215, becomes the new last row, 216) and
# extends the used range to A1:R216.
$ws.Rows(116).Insert()

# Fill in the newly inserted row 116 with the new record's data.
$ws.Range("A116").Value = 10
$ws.Range("B116").Value = "Vega Modelo de Temuco"
$ws.Range("C116").Value = "La Araucanía"
$ws.Range("D116").Value = 44634
$ws.Range("E116").Value = 9
$ws.Range("F116").Value = 100112052
$ws.Range("G116").Value = "Albahaca"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 50
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 5000
$ws.Range("N116").Value = "`$/paquete"
$ws.Range("O116").Value = "Región del Maule"
$ws.Range("P116").Value = 5000
$ws.Range("Q116").Value = 1
$ws.Range("R116").Value = "Hortaliza"
